$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.388.79"
$ws.Range("E2").Value = "  +4.10%  "
$ws.Range("D3").Value = "1.801.46"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'314.87"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.5507"
$ws.Range("E7").Value = "  +4.19%  "
$ws.Range("D8").Value = "'0.3855"
$ws.Range("E8").Value = "  +5.69%  "
$ws.Range("D9").Value = "'0.07600"
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").Value = "'42.53"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("D11").Value = "'1.127"
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "'21.21"
$ws.Range("E13").Value = "  +3.92%  "
$ws.Range("D14").Value = "'6.185"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").Value = "'7.423"
$ws.Range("E15").Value = "  +6.93%  "
$ws.Range("D16").Value = "1.809.99"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "'92.03"
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("D18").Value = "'0.00001072"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("D19").Value = "'0.06443"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "'0.9996"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E21").Value = "  +3.94%  "
$ws.Range("D22").Value = "'5.973"
$ws.Range("D23").Value = "28.432.57"
$ws.Range("E23").Value = "  +3.98%  "
$ws.Range("D24").Value = "'11.45"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").Value = "'2.134"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").Value = "'158.83"
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("D27").Value = "'20.70"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("D28").Value = "'2.412"
$ws.Range("E28").Value = "  +3.72%  "
$ws.Range("D29").Value = "2.016.20"
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("D30").Value = "'123.98"
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("E31").Value = "  +6.03%  "
$ws.Range("D32").Value = "'0.1025"
$ws.Range("E32").Value = "  +5.08%  "
$ws.Range("D33").Value = "'5.774"
$ws.Range("E33").Value = "  +3.87%  "
$ws.Range("D34").Value = "'3.690"
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("D35").Value = "'0.2312"
$ws.Range("E35").Value = "  +14.47%  "
$ws.Range("D36").Value = "'0.06403"
$ws.Range("E36").Value = "  +7.44%  "
$ws.Range("D37").Value = "'0.02327"
$ws.Range("E37").Value = "  +4.46%  "
$ws.Range("D38").Value = "'5.159"
$ws.Range("E38").Value = "  +6.65%  "
$ws.Range("D39").Value = "'8.775"
$ws.Range("E39").Value = "  +9.06%  "
$ws.Range("D40").Value = "'11.64"
$ws.Range("E40").Value = "  +3.95%  "
$ws.Range("D41").Value = "'0.6423"
$ws.Range("E41").Value = "  +4.90%  "
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").Value = "'0.9994"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.159"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("E44").Value = "  -3.12%  "
$ws.Range("D45").Value = "'13.59"
$ws.Range("E45").Value = "  +3.89%  "
$ws.Range("D46").Value = "'0.5981"
$ws.Range("E46").Value = "  +4.12%  "
$ws.Range("D47").Value = "'3.681"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D48").Value = "'126.69"
$ws.Range("E48").Value = "  +4.80%  "
$ws.Range("D49").Value = "'1.985"
$ws.Range("E49").Value = "  +5.09%  "
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("E51").Value = "  +2.83%  "
